$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 02:27"

# Update country rows whose rank/values changed (B..H) due to refreshed stats
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$rowVals = @(4432102, 60263, 2133582, 2148102, 0, 570, 150418)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(4, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(5, 1).Value = "Brasil"
$rowVals = @(2443480, 23579, 1667667, 688134, 0, 627, 87679)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(5, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(23, 1).Value = "Argentina"
$rowVals = @(167416, 4890, 72575, 91782, 0, 120, 3059)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(23, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(24, 1).Value = "Canada"
$rowVals = @(114597, 686, 99860, 5836, 0, 11, 8901)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(24, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(56, 1).Value = "Ghana"
$rowVals = @(33624, 655, 29801, 3655, 0, 0, 168)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(56, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(57, 1).Value = "Kirguistan"
$rowVals = @(33296, 483, 21205, 10790, 0, 24, 1301)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(57, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(64, 1).Value = "Uzbekistan"
$rowVals = @(21209, 678, 11674, 9414, 0, 5, 121)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(64, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(69, 1).Value = "Camerun"
$rowVals = @(17110, 402, 14539, 2180, 0, 6, 391)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(69, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(102, 1).Value = "Paraguay"
$rowVals = @(4548, 104, 2905, 1600, 0, 2, 43)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(102, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(103, 1).Value = "Hungria"
$rowVals = @(4448, 13, 3329, 523, 0, 0, 596)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(103, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(110, 1).Value = "Congo"
$rowVals = @(3200, 162, 829, 2317, 0, 3, 54)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(110, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(111, 1).Value = "Somalia"
$rowVals = @(3196, 18, 1543, 1560, 0, 0, 93)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(111, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(112, 1).Value = "Guinea Ecuatorial"
$rowVals = @(3071, 0, 842, 2178, 0, 0, 51)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(112, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(113, 1).Value = "Mayotte"
$rowVals = @(2900, 38, 2672, 190, 0, 0, 38)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(113, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(114, 1).Value = "Montenegro"
$rowVals = @(2893, 94, 809, 2039, 0, 2, 45)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(114, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(115, 1).Value = "Libia"
$rowVals = @(2827, 158, 577, 2186, 0, 4, 64)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(115, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(116, 1).Value = "Sri Lanka"
$rowVals = @(2805, 23, 2121, 673, 0, 0, 11)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(116, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(118, 1).Value = "Zimbabue"
$rowVals = @(2704, 192, 542, 2126, 0, 2, 36)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(118, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(119, 1).Value = "Cuba"
$rowVals = @(2532, 37, 2351, 94, 0, 0, 87)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(119, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(120, 1).Value = "Mali"
$rowVals = @(2513, 3, 1913, 476, 0, 1, 124)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(120, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(137, 1).Value = "Surinam"
$rowVals = @(1483, 44, 925, 534, 0, 1, 24)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(137, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(138, 1).Value = "Tunez"
$rowVals = @(1455, 3, 1157, 248, 0, 0, 50)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(138, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(140, 1).Value = "Uruguay"
$rowVals = @(1202, 10, 951, 216, 0, 1, 35)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(140, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(147, 1).Value = "Republica de Chipre"
$rowVals = @(1060, 3, 852, 189, 0, 0, 19)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(147, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(151, 1).Value = "Togo"
$rowVals = @(874, 6, 607, 249, 0, 0, 18)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(151, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(152, 1).Value = "Santo Tome y Principe"
$rowVals = @(865, 2, 734, 117, 0, 0, 14)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(152, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(190, 1).Value = "Antigua y Barbuda"
$rowVals = @(86, 4, 65, 18, 0, 0, 3)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(190, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(191, 1).Value = "Liechtenstein"
$rowVals = @(86, 0, 81, 4, 0, 0, 1)
for ($c = 2; $c -le 8; $c++) { $ws.Cells.Item(191, $c).Value = $rowVals[$c - 2] }
